$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 27; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Especial"; M = 300; N = 2400; O = 2500; P = 2450; Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 2450; T = 1 },
    @{ Row = 28; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Primera"; M = 400; N = 2100; O = 2200; P = 2150; Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 2150; T = 1 },
    @{ Row = 29; A = 2; B = "Comercializadora del Agro de Limarí"; C = "Coquimbo"; D = 44461; E = 4; F = "Fruta"; G = 100107; H = "Otros"; I = 100107002; J = "Chirimoya"; K = "Cultivar IV Región"; L = "Segunda"; M = 400; N = 1800; O = 1900; P = 1850; Q = "$/kilo (en caja de 15 kilos)"; R = "Provincia de Limarí"; S = 1850; T = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
}
